$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> @(D_new, E_new). Use $null to skip updating a column.
$updates = @{
    2  = @("310.72", "1.48%")
    3  = @("41.33", "0.77%")
    4  = @("5.137", "0.59%")
    5  = @("0.07699", "1.15%")
    6  = @("1.632", "0.94%")
    7  = @("0.9242", "1.90%")
    9  = @("0.1211", "18.53%")
    10 = @("0.1840", "5.13%")
    11 = @("0.09128", "0.23%")
    12 = @("0.04220", "-1.02%")
    13 = @($null, "-0.57%")
    14 = @("0.001250", "1.79%")
    15 = @("0.005744", "-1.96%")
    17 = @("4.313", "1.15%")
    18 = @("0.3336", "1.90%")
    19 = @("6.939", "6.08%")
    20 = @("0.1390", "2.55%")
    21 = @($null, "-1.81%")
    22 = @("0.04050", "-3.05%")
    23 = @("0.001260", "2.71%")
    24 = @($null, "0.64%")
    25 = @($null, "-2.50%")
    26 = @($null, "24.63%")
    38 = @("0.02469", "3.93%")
    39 = @("0.05278", "2.94%")
    40 = @("0.007835", "0.99%")
    41 = @("0.1315", "1.26%")
    42 = @("0.006797", "-3.03%")
    43 = @($null, "-3.99%")
    44 = @("0.008189", "-3.09%")
    45 = @("0.3100", "-6.91%")
    46 = @("0.00006725", "5.80%")
    47 = @($null, "-0.19%")
    48 = @("0.1699", "515.16%")
    49 = @("0.004096", "-6.98%")
    50 = @("0.00002098", "-0.19%")
    51 = @("0.0001998", "-0.19%")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($null -ne $dVal) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
    }
    if ($null -ne $eVal) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $eVal
    }
}
